# inlineForeignTabs 可编辑聚合表格, isDeleteCascade 级联删除
#
# The header row of this import-template sheet lists one cell per
# generated column, each cell holding an EJS comment/validation snippet
# (via a shared string). Drop the "is_locked" (D1) and "is_enabled" (E1)
# template columns entirely: the cells to their right ("order_by", "rem")
# shift left into D1/E1, and the now-unreferenced shared strings for
# is_locked/is_enabled are dropped when the workbook is re-saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1:E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
